$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hungary NB I")

# Copy row 187 formatting down to the new rows (188-193) before filling values
$ws.Range("A187:AB187").Copy() | Out-Null
$ws.Range("A188:AB193").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 188
$ws.Range("A188").Value = 186
$ws.Range("B188").Value = 7021685
$ws.Range("C188").Value = "Hungary NB I"
$ws.Range("D188").Value = 45422.63541666666
$ws.Range("E188").Value = "Ujpest"
$ws.Range("F188").Value = "Kecskemeti TE"
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 3
$ws.Range("I188").Value = "A"
$ws.Range("J188").Value = 2.2
$ws.Range("K188").Value = 3.4
$ws.Range("L188").Value = 3.2
$ws.Range("M188").Value = 2.6
$ws.Range("N188").Value = 3.4
$ws.Range("O188").Value = 2.625
$ws.Range("P188").Value = 0
$ws.Range("Q188").Value = 1.925
$ws.Range("R188").Value = 1.925
$ws.Range("S188").Value = 2.5
$ws.Range("T188").Value = 1.825
$ws.Range("U188").Value = 2.025
$ws.Range("V188").Value = -1
$ws.Range("W188").Value = -1
$ws.Range("X188").Value = 1.625
$ws.Range("Y188").Value = -1
$ws.Range("Z188").Value = 0.925
$ws.Range("AA188").Value = 0.825
$ws.Range("AB188").Value = -1

# Row 189
$ws.Range("A189").Value = 187
$ws.Range("B189").Value = 7021686
$ws.Range("C189").Value = "Hungary NB I"
$ws.Range("D189").Value = 45423.4375
$ws.Range("E189").Value = "MTK Budapest"
$ws.Range("F189").Value = "Puskas Academy"
$ws.Range("G189").Value = 1
$ws.Range("H189").Value = 3
$ws.Range("I189").Value = "A"
$ws.Range("J189").Value = 3.6
$ws.Range("K189").Value = 4
$ws.Range("L189").Value = 1.85
$ws.Range("M189").Value = 5.75
$ws.Range("N189").Value = 5
$ws.Range("O189").Value = 1.444
$ws.Range("P189").Value = 1.25
$ws.Range("Q189").Value = 1.85
$ws.Range("R189").Value = 2
$ws.Range("S189").Value = 3
$ws.Range("T189").Value = 1.925
$ws.Range("U189").Value = 1.925
$ws.Range("V189").Value = -1
$ws.Range("W189").Value = -1
$ws.Range("X189").Value = 0.444
$ws.Range("Y189").Value = -1
$ws.Range("Z189").Value = 1
$ws.Range("AA189").Value = 0.925
$ws.Range("AB189").Value = -1

# Row 190
$ws.Range("A190").Value = 188
$ws.Range("B190").Value = 7028360
$ws.Range("C190").Value = "Hungary NB I"
$ws.Range("D190").Value = 45423.52083333334
$ws.Range("E190").Value = "Zalaegerszegi TE"
$ws.Range("F190").Value = "Paksi"
$ws.Range("G190").Value = 1
$ws.Range("H190").Value = 1
$ws.Range("I190").Value = "D"
$ws.Range("J190").Value = 2.6
$ws.Range("K190").Value = 3.4
$ws.Range("L190").Value = 2.6
$ws.Range("M190").Value = 3.1
$ws.Range("N190").Value = 3.5
$ws.Range("O190").Value = 2.2
$ws.Range("P190").Value = 0.25
$ws.Range("Q190").Value = 1.925
$ws.Range("R190").Value = 1.925
$ws.Range("S190").Value = 3
$ws.Range("T190").Value = 1.975
$ws.Range("U190").Value = 1.875
$ws.Range("V190").Value = -1
$ws.Range("W190").Value = 2.5
$ws.Range("X190").Value = -1
$ws.Range("Y190").Value = 0.4625
$ws.Range("Z190").Value = -0.5
$ws.Range("AA190").Value = -1
$ws.Range("AB190").Value = 0.875

# Row 191
$ws.Range("A191").Value = 189
$ws.Range("B191").Value = 7024477
$ws.Range("C191").Value = "Hungary NB I"
$ws.Range("D191").Value = 45423.60416666666
$ws.Range("E191").Value = "Diosgyori VTK"
$ws.Range("F191").Value = "Ferencvarosi TC"
$ws.Range("G191").Value = 2
$ws.Range("H191").Value = 0
$ws.Range("I191").Value = "H"
$ws.Range("J191").Value = 4.5
$ws.Range("K191").Value = 4.2
$ws.Range("L191").Value = 1.65
$ws.Range("M191").Value = 4.333
$ws.Range("N191").Value = 4.75
$ws.Range("O191").Value = 1.6
$ws.Range("P191").Value = 1
$ws.Range("Q191").Value = 1.85
$ws.Range("R191").Value = 2
$ws.Range("S191").Value = 3
$ws.Range("T191").Value = 1.85
$ws.Range("U191").Value = 2
$ws.Range("V191").Value = 3.333
$ws.Range("W191").Value = -1
$ws.Range("X191").Value = -1
$ws.Range("Y191").Value = 0.8500000000000001
$ws.Range("Z191").Value = -1
$ws.Range("AA191").Value = -1
$ws.Range("AB191").Value = 1

# Row 192
$ws.Range("A192").Value = 190
$ws.Range("B192").Value = 7024498
$ws.Range("C192").Value = "Hungary NB I"
$ws.Range("D192").Value = 45424.38541666666
$ws.Range("E192").Value = "Kisvarda FC"
$ws.Range("F192").Value = "Mezokovesd Zsory"
$ws.Range("G192").Value = 4
$ws.Range("H192").Value = 3
$ws.Range("I192").Value = "H"
$ws.Range("J192").Value = 1.95
$ws.Range("K192").Value = 3.9
$ws.Range("L192").Value = 3.4
$ws.Range("M192").Value = 1.6
$ws.Range("N192").Value = 4.5
$ws.Range("O192").Value = 4.75
$ws.Range("P192").Value = -1
$ws.Range("Q192").Value = 2.05
$ws.Range("R192").Value = 1.8
$ws.Range("S192").Value = 3
$ws.Range("T192").Value = 1.85
$ws.Range("U192").Value = 2
$ws.Range("V192").Value = 0.6000000000000001
$ws.Range("W192").Value = -1
$ws.Range("X192").Value = -1
$ws.Range("Y192").Value = 0
$ws.Range("Z192").Value = 0
$ws.Range("AA192").Value = 0.8500000000000001
$ws.Range("AB192").Value = -1

# Row 193
$ws.Range("A193").Value = 191
$ws.Range("B193").Value = 7028515
$ws.Range("C193").Value = "Hungary NB I"
$ws.Range("D193").Value = 45424.60416666666
$ws.Range("E193").Value = "Debreceni VSC"
$ws.Range("F193").Value = "MOL Fehervar FC"
$ws.Range("G193").Value = 1
$ws.Range("H193").Value = 0
$ws.Range("I193").Value = "H"
$ws.Range("J193").Value = 2.6
$ws.Range("K193").Value = 3.4
$ws.Range("L193").Value = 2.6
$ws.Range("M193").Value = 2
$ws.Range("N193").Value = 3.5
$ws.Range("O193").Value = 3.25
$ws.Range("P193").Value = -0.25
$ws.Range("Q193").Value = 1.825
$ws.Range("R193").Value = 2.025
$ws.Range("S193").Value = 2.75
$ws.Range("T193").Value = 1.85
$ws.Range("U193").Value = 2
$ws.Range("V193").Value = 1
$ws.Range("W193").Value = -1
$ws.Range("X193").Value = -1
$ws.Range("Y193").Value = 0.825
$ws.Range("Z193").Value = -1
$ws.Range("AA193").Value = -1
$ws.Range("AB193").Value = 1
